$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3325
$ws1.Range("F6").Value = 2288
$ws1.Range("F8").Value = 90
$ws1.Range("F9").Value = 28
$ws1.Range("F11").Value = 75
$ws1.Range("F14").Value = 442
$ws1.Range("F15").Value = 41

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3325
$ws4.Range("F7").Value = 2288
$ws4.Range("F9").Value = 90
$ws4.Range("F10").Value = 28
$ws4.Range("F12").Value = 75
$ws4.Range("F15").Value = 442
$ws4.Range("F16").Value = 41
